$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 75.7
$ws.Range("C3").Value = 75.90000000000001
$ws.Range("C4").Value = 73.90000000000001
$ws.Range("C5").Value = 80.2
$ws.Range("C6").Value = 80.2
$ws.Range("C7").Value = 80.5
